# Update "Forecast Comparison" sheet with a new Week_Start_Date column and
# corrected, shorter week labels (W01 -> W1, etc.), per "Update with Correct
# Forecast output".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the existing ASIN column (B) to hold the week
# start date. Everything from the old B column onward shifts right by one.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week label (A) and corresponding ISO week-start date (B) for each of the
# 16 forecast rows (rows 2-17).
$weeks = @("W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8", "W9", "W10", "W11", "W12", "W13", "W14", "W15", "W16")
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

# Format the new column as Text first so the "yyyy-mm-dd" strings are stored
# as literal text (matching the source data) instead of being auto-converted
# into date serial numbers.
$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $weeks[$i]
    $ws.Cells.Item($r, 2).Value = $weekStartDates[$i]
}

# is_holiday_week (now column J after the insert) should be a boolean FALSE,
# not the numeric 0 it held before.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 10).Value = $false
}
